# Update "想去人数" (interest count) figures in the 展览 (sheet 1) and
# 全部类型 (sheet 4) worksheets to match the latest scrape snapshot.
$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("F2").Value = 334
    $ws.Range("F4").Value = 1510
    $ws.Range("F9").Value = 320
}
